$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.36664133333333
$ws.Range("H2").Value = 70.099924
$ws.Range("I2").Value = 0.7520698922374942
$ws.Range("J2").Value = 0.7520698922374941
$ws.Range("M2").Value = 8.226724333333333
$ws.Range("N2").Value = 24.680173
$ws.Range("O2").Value = 0.06198126651953669
$ws.Range("P2").Value = 0.06198126651953669
$ws.Range("Q2").Value = 192.2309168452058
$ws.Range("R2").Value = 1730.078251606852
$ws.Range("S2").Value = 0.04661424443209137
$ws.Range("T2").Value = 0.04661424443209136

$ws.Range("G3").Value = 23.36664133333333
$ws.Range("H3").Value = 70.099924
$ws.Range("I3").Value = 0.7520698922374942
$ws.Range("J3").Value = 0.7520698922374941
$ws.Range("O3").Value = 0.6623065855236785
$ws.Range("P3").Value = 0.6623065855236785
$ws.Range("Q3").Value = 2054.101332822944
$ws.Range("R3").Value = 18486.91199540649
$ws.Range("S3").Value = 0.4981008424029756
$ws.Range("T3").Value = 0.4981008424029755

$ws.Range("G4").Value = 23.36664133333333
$ws.Range("H4").Value = 70.099924
$ws.Range("I4").Value = 0.7520698922374942
$ws.Range("J4").Value = 0.7520698922374941
$ws.Range("M4").Value = 36.43008433333333
$ws.Range("N4").Value = 109.290253
$ws.Range("O4").Value = 0.2744692388979848
$ws.Range("P4").Value = 0.2744692388979848
$ws.Range("Q4").Value = 851.2487143600858
$ws.Range("R4").Value = 7661.238429240772
$ws.Range("S4").Value = 0.2064200509205145
$ws.Range("T4").Value = 0.2064200509205145

$ws.Range("G5").Value = 23.36664133333333
$ws.Range("H5").Value = 70.099924
$ws.Range("I5").Value = 0.7520698922374942
$ws.Range("J5").Value = 0.7520698922374941
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1649703333333333
$ws.Range("N5").Value = 0.494911
$ws.Range("O5").Value = 0.00124290905879997
$ws.Range("P5").Value = 0.00124290905879997
$ws.Range("Q5").Value = 3.854802609640445
$ws.Range("R5").Value = 34.693223486764
$ws.Range("S5").Value = 0.0009347544819126985
$ws.Range("T5").Value = 0.0009347544819126983

$ws.Range("I6").Value = 0.2155695244255668
$ws.Range("J6").Value = 0.2155695244255668
$ws.Range("M6").Value = 8.226724333333333
$ws.Range("N6").Value = 24.680173
$ws.Range("O6").Value = 0.06198126651953669
$ws.Range("P6").Value = 0.06198126651953669
$ws.Range("Q6").Value = 55.10010140271078
$ws.Range("R6").Value = 495.900912624397
$ws.Range("S6").Value = 0.01336127214691083
$ws.Range("T6").Value = 0.01336127214691083

$ws.Range("I7").Value = 0.2155695244255668
$ws.Range("J7").Value = 0.2155695244255668
$ws.Range("O7").Value = 0.6623065855236785
$ws.Range("P7").Value = 0.6623065855236785
$ws.Range("S7").Value = 0.1427731156652604
$ws.Range("T7").Value = 0.1427731156652604

$ws.Range("I8").Value = 0.2155695244255668
$ws.Range("J8").Value = 0.2155695244255668
$ws.Range("M8").Value = 36.43008433333333
$ws.Range("N8").Value = 109.290253
$ws.Range("O8").Value = 0.2744692388979848
$ws.Range("P8").Value = 0.2744692388979848
$ws.Range("Q8").Value = 243.9976422623908
$ws.Range("R8").Value = 2195.978780361517
$ws.Range("S8").Value = 0.05916720329868587
$ws.Range("T8").Value = 0.05916720329868587

$ws.Range("I9").Value = 0.2155695244255668
$ws.Range("J9").Value = 0.2155695244255668
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1649703333333333
$ws.Range("N9").Value = 0.494911
$ws.Range("O9").Value = 0.00124290905879997
$ws.Range("P9").Value = 0.00124290905879997
$ws.Range("Q9").Value = 1.104921196675444
$ws.Range("R9").Value = 9.944290770079
$ws.Range("S9").Value = 0.0002679333147097383
$ws.Range("T9").Value = 0.0002679333147097383

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.005436
$ws.Range("H10").Value = 3.016308
$ws.Range("I10").Value = 0.03236058333693902
$ws.Range("J10").Value = 0.03236058333693902
$ws.Range("M10").Value = 8.226724333333333
$ws.Range("N10").Value = 24.680173
$ws.Range("O10").Value = 0.06198126651953669
$ws.Range("P10").Value = 0.06198126651953669
$ws.Range("Q10").Value = 8.271444806809333
$ws.Range("R10").Value = 74.44300326128401
$ws.Range("S10").Value = 0.002005749940534495
$ws.Range("T10").Value = 0.002005749940534495

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.005436
$ws.Range("H11").Value = 3.016308
$ws.Range("I11").Value = 0.03236058333693902
$ws.Range("J11").Value = 0.03236058333693902
$ws.Range("O11").Value = 0.6623065855236785
$ws.Range("P11").Value = 0.6623065855236785
$ws.Range("Q11").Value = 88.38529244346266
$ws.Range("R11").Value = 795.467631991164
$ws.Range("S11").Value = 0.02143262745544253
$ws.Range("T11").Value = 0.02143262745544253

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.005436
$ws.Range("H12").Value = 3.016308
$ws.Range("I12").Value = 0.03236058333693902
$ws.Range("J12").Value = 0.03236058333693902
$ws.Range("M12").Value = 36.43008433333333
$ws.Range("N12").Value = 109.290253
$ws.Range("O12").Value = 0.2744692388979848
$ws.Range("P12").Value = 0.2744692388979848
$ws.Range("Q12").Value = 36.62811827176933
$ws.Range("R12").Value = 329.653064445924
$ws.Range("S12").Value = 0.008881984678784462
$ws.Range("T12").Value = 0.008881984678784462

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.005436
$ws.Range("H13").Value = 3.016308
$ws.Range("I13").Value = 0.03236058333693902
$ws.Range("J13").Value = 0.03236058333693902
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1649703333333333
$ws.Range("N13").Value = 0.494911
$ws.Range("O13").Value = 0.00124290905879997
$ws.Range("P13").Value = 0.00124290905879997
$ws.Range("Q13").Value = 0.1658671120653333
$ws.Range("R13").Value = 1.492804008588
$ws.Range("S13").Value = 0.00004022126217753285
$ws.Range("T13").Value = 0.00004022126217753285
